# Auto update Excel log
# - Append new sensor/alert rows to ALERTS, Proximity and mmWave sheets.
# - Add a new "Camera" sheet with its own header + first capture row.

$wb = $excel.ActiveWorkbook

function Add-LogRow {
    param($ws, [int]$row, [string]$date, [string]$timestamp, [string]$hour, [string]$location, [string]$value, [string]$status)

    # Column A holds plain "YYYY-MM-DD" text in this log. Excel's COM layer
    # auto-coerces that pattern into a real date serial on assignment, so
    # force the cell to text first to preserve the literal string.
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $date

    $ws.Cells.Item($row, 2).Value = $timestamp
    $ws.Cells.Item($row, 3).Value = $hour
    $ws.Cells.Item($row, 4).Value = $location
    $ws.Cells.Item($row, 5).Value = $value
    $ws.Cells.Item($row, 6).Value = $status
}

# ---------------------------------------------------------------------------
# ALERTS: new WARNING row
# ---------------------------------------------------------------------------
$alerts = $wb.Worksheets.Item("ALERTS")
Add-LogRow $alerts 5 "2026-01-28" "17:36:33" "17:00" "Bathroom" "WARNING" "Bathroom Humidity > 90.0% for 20s with NO MOTION. Alerting."

# ---------------------------------------------------------------------------
# Proximity: door enter/exit events
# ---------------------------------------------------------------------------
$proximity = $wb.Worksheets.Item("Proximity")
Add-LogRow $proximity 3 "2026-01-28" "17:37:22" "17:00" "Bathroom Door" "EXIT" "User EXITED Bathroom"
Add-LogRow $proximity 4 "2026-01-28" "17:37:27" "17:00" "Living Room Main Door" "ENTER" "User ENTERED Living Room Main Door"
Add-LogRow $proximity 5 "2026-01-28" "17:37:29" "17:00" "Living Room Main Door" "EXIT" "User EXITED Living Room Main Door"

# ---------------------------------------------------------------------------
# mmWave: presence timeline continuation
# ---------------------------------------------------------------------------
$mmwave = $wb.Worksheets.Item("mmWave")

$mmwaveRows = @(
    @("17:36:33", "PRESENCE", "Active"),
    @("17:36:33", "PRESENCE", "Active"),
    @("17:36:33", "PRESENCE", "Active"),
    @("17:36:36", "PRESENCE", "Active"),
    @("17:36:39", "PRESENCE", "Active"),
    @("17:36:42", "PRESENCE", "Active"),
    @("17:36:45", "PRESENCE", "Active"),
    @("17:36:48", "PRESENCE", "Active"),
    @("17:36:51", "PRESENCE", "Active"),
    @("17:36:54", "PRESENCE", "Active"),
    @("17:36:57", "PRESENCE", "Active"),
    @("17:37:00", "PRESENCE", "Active"),
    @("17:37:03", "PRESENCE", "Active"),
    @("17:37:06", "PRESENCE", "Active"),
    @("17:37:09", "PRESENCE", "Active"),
    @("17:37:12", "PRESENCE", "Active"),
    @("17:37:15", "PRESENCE", "Active"),
    @("17:37:18", "PRESENCE", "Active"),
    @("17:37:21", "NO_PRESENCE", "Inactive"),
    @("17:37:24", "NO_PRESENCE", "Inactive"),
    @("17:37:27", "NO_PRESENCE", "Inactive"),
    @("17:37:30", "NO_PRESENCE", "Inactive")
)

$r = 56
foreach ($row in $mmwaveRows) {
    Add-LogRow $mmwave $r "2026-01-28" $row[0] "17:00" "Living Room" $row[1] $row[2]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Camera: brand new sheet, appended after mmWave (last existing tab)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$camera = $wb.Worksheets.Add($null, $lastSheet)
$camera.Name = "Camera"

$camera.Cells.Item(1, 1).Value = "Date"
$camera.Cells.Item(1, 2).Value = "Timestamp"
$camera.Cells.Item(1, 3).Value = "Hour"
$camera.Cells.Item(1, 4).Value = "Location"
$camera.Cells.Item(1, 5).Value = "Value"
$camera.Cells.Item(1, 6).Value = "Status"

Add-LogRow $camera 2 "2026-01-28" "17:37:29" "17:00" "Living Room Main Door" "Image Captured" "Active"
